$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FormatFrom($srcAddr, $dstAddr) {
  $ws.Range($srcAddr).Copy()
  $ws.Range($dstAddr).PasteSpecial(-4122)
  $excel.CutCopyMode = $false
}

# ============ Row 6 block: D6 -> grey, E6 removed ============
$ws.Range("E6:E11").UnMerge()
$ws.Range("E6:E11").ClearContents()
$ws.Range("E6:E11").ClearFormats()
Set-FormatFrom "E14" "D6"
$ws.Range("D6").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n09:00-10:30`nroom:Grey classroom: 203"

# ============ Row 14 block: add B14,C14,D14 (merge FIRST, then format) ============
$ws.Range("B14:B19").Merge()
$ws.Range("C14:C19").Merge()
$ws.Range("D14:D19").Merge()
Set-FormatFrom "E14" "B14"
$ws.Range("B14").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n11:00-12:30`nroom:Grey classroom: 203"
Set-FormatFrom "C6" "C14"
$ws.Range("C14").Value = "Precalculus`n`n11:00-12:30`nroom:Green classroom: 204"
Set-FormatFrom "E14" "D14"
$ws.Range("D14").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n11:00-12:30`nroom:Grey classroom: 203"

# ---- rows 15-18 (top border) and 19 (bottom border) for B, C, D ----
foreach ($r in 15..18) {
  Set-FormatFrom "B7" "B$r"
  Set-FormatFrom "B7" "C$r"
  Set-FormatFrom "B7" "D$r"
}
Set-FormatFrom "B11" "B19"
Set-FormatFrom "B11" "C19"
Set-FormatFrom "B11" "D19"

# ============ Row 26 block: add B26 (green), D26 -> grey, remove E26 ============
$ws.Range("E26:E31").UnMerge()
$ws.Range("E26:E31").ClearContents()
$ws.Range("E26:E31").ClearFormats()
$ws.Range("B26:B31").Merge()
Set-FormatFrom "C6" "B26"
$ws.Range("B26").Value = "Precalculus`n`n14:00-15:30`nroom:Green classroom: 204"
Set-FormatFrom "E14" "D26"
$ws.Range("D26").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n14:00-15:30`nroom:Grey classroom: 203"

# ---- rows 27-30 (top border) and 31 (bottom border) for B ----
foreach ($r in 27..30) {
  Set-FormatFrom "B7" "B$r"
}
Set-FormatFrom "B11" "B31"

# ============ Row 34 block: D34 -> red Scientific Inquiry, E34 -> red Academic Writing ============
Set-FormatFrom "F6" "D34"
$ws.Range("D34").Value = "Scientific Inquiry: Beyond the Visible`n`n16:00-17:30`nroom:Red classroom: 201"
Set-FormatFrom "F6" "E34"
$ws.Range("E34").Value = "Academic Writing: Research, Fiction and Nonfiction`n`n16:00-17:30`nroom:Red classroom: 201"

Write-Host "Done"
